# D1 - Se repérer dans le temps (A trou)
# "Ajout chapitre 12-3 6ème"
#
# 1) The section I title changes from "Système décimal" to "Durée".
# 2) The worked example sentence is corrected from the singular
#    "ces deux instant" to the plural "ces deux instants".

$d = $word.ActiveDocument

# --- 1. Title: "Système décimal" -> "Durée" -----------------------------
$null = $d.Content.Find.Execute("Système décimal", $false, $false, $false, $false, $false, $true, 1, $false, "Durée", 2)

# --- 2. "deux instant est" -> "deux instants est" ------------------------
$null = $d.Content.Find.Execute("deux instant est", $false, $false, $false, $false, $false, $true, 1, $false, "deux instants est", 2)
